$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Activate()

# Widen column C
$ws.Columns.Item(3).ColumnWidth = 37.20703125

# Update the two numeric values in column C
$ws.Range("C2").Value = 0.00034101243963859801
$ws.Range("C3").Value = 34.439998626708899

# Make row 3 taller
$ws.Rows.Item(3).RowHeight = 17.7

# Apply a Verdana/14/black font to C3 via a throwaway named style so the
# font table only grows by the one font that's actually needed
$c3 = $ws.Range("C3")
$style = $wb.Styles.Add("TempFontStyle")
$style.Font.Name = "Verdana"
$style.Font.Size = 14
$style.Font.Color = 0
$c3.Style = "TempFontStyle"
$wb.Styles.Item("TempFontStyle").Delete()

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1

# Move the selection
$ws.Range("I20").Select()
